# Insert a new data row at row 82 (pushing the existing rows 82-163 down to
# 83-164) and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(82).Insert()

$ws.Cells.Item(82, 1).Value  = 8
$ws.Cells.Item(82, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(82, 3).Value  = "Coquimbo"
$ws.Cells.Item(82, 4).Value  = 44827
$ws.Cells.Item(82, 5).Value  = 4
$ws.Cells.Item(82, 6).Value  = 100112001
$ws.Cells.Item(82, 7).Value  = "Berenjena"
$ws.Cells.Item(82, 8).Value  = "Sin especificar"
$ws.Cells.Item(82, 9).Value  = "Primera"
$ws.Cells.Item(82, 10).Value = 400
$ws.Cells.Item(82, 11).Value = 10500
$ws.Cells.Item(82, 12).Value = 11000
$ws.Cells.Item(82, 13).Value = 10750
$ws.Cells.Item(82, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(82, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(82, 16).Value = 269
$ws.Cells.Item(82, 17).Value = 40
$ws.Cells.Item(82, 18).Value = "Hortaliza"
